$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The supervisor "Gabby James" has taken over as supervisor/contact for all
# employees in the sheet, so every row's SupervisorEmail value becomes the
# same address, and the per-row hyperlinks collapse into a single hyperlink
# for C2 and one combined hyperlink covering C3:C21.

# Replace all SupervisorEmail cell values (column C, rows 2-21) with the
# new, single supervisor email address.
$ws.Range("C2:C21").Value = "gjames@bhnstl.org"

# Remove all of the old, per-row mailto hyperlinks before adding the new
# consolidated ones.
$ws.Range("A1:E21").Hyperlinks.Delete()

# Re-create the hyperlinks: C2 keeps its own hyperlink, while C3:C21 share
# a single hyperlink definition, matching the target layout.
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:gjames@bhnstl.org", "", "", "gjames@bhnstl.org")
$ws.Hyperlinks.Add($ws.Range("C3:C21"), "mailto:gjames@bhnstl.org", "", "", "gjames@bhnstl.org")

# Adding hyperlinks can introduce a new "Hyperlink" flavored style for the
# first cell of each range; reapply the Hyperlink style uniformly across
# the whole column so every cell keeps the original formatting/style index.
$ws.Range("C2:C21").Style = "Hyperlink"

# Reflect the final cell selection recorded in the saved workbook.
$ws.Range("G20").Select()
